# The document contains a single 46-row, 1-column table. Benchmark numbers
# were re-derived, shifting several values and collapsing the last three
# multi-run rows (which held whole tab-separated summary lines) down to the
# single values that used to live in rows 1-3.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "978"
$t.Cell(5, 1).Range.Text  = "0.00001"
$t.Cell(6, 1).Range.Text  = "0.00256"
$t.Cell(7, 1).Range.Text  = "0.00013"
$t.Cell(8, 1).Range.Text  = "0.00009"
$t.Cell(9, 1).Range.Text  = "0.00021"
$t.Cell(10, 1).Range.Text = "0.00024"
$t.Cell(11, 1).Range.Text = "0.00037"
$t.Cell(12, 1).Range.Text = "0.15467"

$t.Cell(44, 1).Range.Text = "99.84"
$t.Cell(45, 1).Range.Text = "0.15"
$t.Cell(46, 1).Range.Text = "98"
